$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141, shifting existing rows 141:151 down to 142:152
$ws.Rows("141").Insert()

# Populate the new row 141 with the latest week's data for this market/product
$ws.Range("A141").Value = 4
$ws.Range("B141").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C141").Value = "Los Lagos"
$ws.Range("D141").Value = 44516
$ws.Range("E141").Value = 10
$ws.Range("F141").Value = 100112024
$ws.Range("G141").Value = "Choclo"
$ws.Range("H141").Value = "Dulce o Americano"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 200
$ws.Range("K141").Value = 26000
$ws.Range("L141").Value = 26000
$ws.Range("M141").Value = 26000
$ws.Range("N141").Value = "$/malla 70 unidades"
$ws.Range("O141").Value = "Región de Arica y Parinacota"
$ws.Range("P141").Value = 371
$ws.Range("Q141").Value = 70
$ws.Range("R141").Value = "Hortaliza"
